# "Add files via upload" — extends the daily price/volume history on the
# "gUSD 26.06.25" sheet (rows 93-103) with 11 more days of data, and moves
# the sheet's selection to reflect the new bottom of the filled range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gUSD 26.06.25")

# New daily rows: Row, C, D, E, F, G
$newRows = @(
    @(93,  90.405299999999997, 36.2363,             8.73, 21.31, 16.989999999999998),
    @(94,  90.409899999999993, 37.698500000000003,  8.85, 22.04, 11.14),
    @(95,  89.385999999999996, 38.912199999999999,  8.86, 21.86, 9.16),
    @(96,  88.315600000000003, 39.662300000000002,  8.8800000000000008, 22.19, 5.59),
    @(97,  88.866900000000001, 41.546300000000002,  9.06, 21.67, 14.57),
    @(98,  83.053799999999995, 42.134399999999999,  8.56, 12.41, 4.2699999999999996),
    @(99,  76.781599999999997, 42.8001,              8,    9.42, 4.74),
    @(100, 75.685699999999997, 44.094999999999999,  7.99, 8.4700000000000006, 10.210000000000001),
    @(101, 72.999600000000001, 45.314700000000002,  7.82, 8.1300000000000008, 8.2100000000000009),
    @(102, 72.515000000000001, 46.799700000000001,  7.88, 8.4700000000000006, 12.29),
    @(103, 71.372600000000006, 47.039299999999997,  7.88, 7.88, 1.41)
)

foreach ($row in $newRows) {
    $r = $row[0]

    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]

    # I = C / $D$3, matching the existing shared formula in I65:I92
    $ws.Range("I$r").Formula = "=C$r/`$D`$3"
    $ws.Range("I$r").NumberFormat = "0.0000"

    # M = C + D, matching the existing shared formula in M65:M92
    $ws.Range("M$r").Formula = "=C$r+D$r"
}

# Reflect the author's new cursor position on the sheet (was C93, now C104 -
# the first still-empty row right after the newly filled data).
[void]$ws.Activate()
[void]$ws.Range("C104").Select()
